$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DE {
    param($Row, $D, $ForceText, $E)

    if ($D -ne $null) {
        $dCell = $ws.Cells.Item($Row, 4)
        if ($ForceText -eq 1) {
            # Force the cell to keep a numeric-looking string as literal text,
            # matching the original file's plain (unstyled) text cells.
            $dCell.NumberFormat = "@"
            $dCell.Value = $D
            $dCell.Style = "Normal"
        } else {
            $dCell.Value = $D
        }
    }

    $ws.Cells.Item($Row, 5).Value = "  $E  "
}

Set-DE 2  "53.578.43"   0 "-4.39%"
Set-DE 3  "2.196.33"    0 "-7.09%"
Set-DE 4  $null         0 "+0.01%"
Set-DE 5  "485.71"      1 "-3.25%"
Set-DE 6  "125.12"      1 "-3.51%"
Set-DE 7  $null         0 "-0.29%"
Set-DE 8  "0.521"       1 "-4.44%"
Set-DE 9  "2.215.55"    0 "-6.40%"
Set-DE 10 "0.0919"      1 "-6.44%"
Set-DE 11 $null         0 "-1.32%"
Set-DE 12 $null         0 "-3.86%"
Set-DE 13 $null         0 "-3.18%"
Set-DE 14 "2.586.23"    0 "-7.12%"
Set-DE 15 "21.12"       1 "-1.43%"
Set-DE 16 "53.504.85"   0 "-4.49%"
Set-DE 17 $null         0 "-3.09%"
Set-DE 18 "2.202.84"    0 "-5.81%"
Set-DE 19 $null         0 "-4.48%"
Set-DE 20 $null         0 "-1.80%"
Set-DE 21 "294.45"      1 "-4.23%"
Set-DE 22 "6.08"        1 "-3.27%"
Set-DE 23 $null         0 "-0.44%"
Set-DE 24 "62.53"       1 "-5.07%"
Set-DE 25 "0.994"       1 "-0.47%"
Set-DE 26 $null         0 "-1.34%"
Set-DE 27 $null         0 "-0.84%"
Set-DE 28 "2.298.86"    0 "-6.94%"
Set-DE 29 "7.01"        1 "-3.07%"
Set-DE 30 "165.03"      1 "-4.17%"
Set-DE 31 $null         0 "-3.75%"
Set-DE 32 $null         0 "-0.19%"
Set-DE 33 "0.0₃0664"    0 "-6.46%"
Set-DE 34 $null         0 "-0.36%"
Set-DE 35 $null         0 "-1.62%"
Set-DE 37 $null         0 "-1.77%"
Set-DE 38 $null         0 "-1.79%"
Set-DE 39 "0.826"       1 "+3.36%"
Set-DE 40 "3.54"        1 "-4.98%"
Set-DE 41 "35.75"       1 "-1.13%"
Set-DE 42 $null         0 "-0.80%"
Set-DE 43 $null         0 "-1.62%"
Set-DE 44 "3.27"        1 "-2.76%"
Set-DE 45 "124.77"      1 "-3.30%"
Set-DE 46 $null         0 "+1.56%"
Set-DE 47 "0.0876"      1 "-2.83%"
Set-DE 48 $null         0 "-5.14%"
Set-DE 49 "232.05"      1 "-2.80%"
Set-DE 50 "0.0470"      1 "-2.53%"
Set-DE 51 $null         0 "-3.12%"
